$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $find
    }
}

# 1. "Ziel" paragraph: fix hyphenation "Georg Simon Ohm-Schule" -> "Georg-Simon-Ohm-Schule"
Replace-Text "Vorbereitungsräumen der Georg Simon Ohm-Schule" "Vorbereitungsräumen der Georg-Simon-Ohm-Schule"

# 2. "Umsetzung" paragraph: "Georg-Simon-Ohm-Berufsschule" -> "Georg-Simon-Ohm-Schule"
Replace-Text "Für die Georg-Simon-Ohm-Berufsschule wurden" "Für die Georg-Simon-Ohm-Schule wurden"

# 3. remove comma after "Smartphones"
Replace-Text "Zudem ist auch das Telefonieren mit Smartphones, über eine kostenlose App möglich." "Zudem ist auch das Telefonieren mit Smartphones über eine kostenlose App möglich."

# 4. add "sind"
Replace-Text "nicht möglich und somit die Gespräche sicher" "nicht möglich und somit sind die Gespräche sicher"

# 5. Anrufbeantworter sentence rewrite
Replace-Text "Diese bestehen aus der Möglichkeit eine individuelle Anrufbeantworter Sprachnachricht aufzunehmen und bei besetzter Leitung eine entsprechende Mail an den nicht erreichten Lehrer mit der aufgenommenen Nachricht per Mail zu senden." "Diese bestehen erstens aus der Möglichkeit eine individuelle Anrufbeantworter Ansage aufzunehmen und zweitens bei besetzter Leitung oder einen verpassten Anruf eine entsprechende Mail an den nicht erreichten Lehrer mit der aufgenommenen Nachricht zu senden."

# 6. Telefonnummer(n) anhäufen -> ansammelt, remove extra comma
Replace-Text "Da sich die Anzahl an Telefonnummer schnell anhäufen, ist ein Telefonbuch, mit allen Nummern vorhanden" "Da sich die Anzahl an Telefonnummern schnell ansammelt, ist ein Telefonbuch mit allen Nummern vorhanden"

# 7. "Zeitlicher Rahmen" paragraph: "3 Teilen" -> "drei Teile"
Replace-Text "Das Projekt gliedert sich in 3 Teilen bezogen" "Das Projekt gliedert sich in drei Teile bezogen"

# 8. "Georg Simon Ohm" -> "Georg-Simon-Ohm-Schule"
Replace-Text "Kompatibel mit den in der Georg Simon Ohm vorhanden Serverinfrastruktur." "Kompatibel mit den in der Georg-Simon-Ohm-Schule vorhanden Serverinfrastruktur."

# 9. "können Experten" -> "wird ein Experte"
Replace-Text "Bei Nachfragen können Experten für einen geringen Aufschlag" "Bei Nachfragen wird ein Experte für einen geringen Aufschlag"
